# Switch the Tubes.xlsx sample data to a single custom row, pointing at a
# different "tube" entry (date 1/19 instead of 5/10), dropping the
# remaining sample rows (2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded sample rows 2 through 10 entirely (shifts cells
# up / shrinks the sheet dimension), keeping row 1 untouched.
$ws.Range("A2:A10").EntireRow.Delete() | Out-Null

# Update the remaining row's date-like text value.
$ws.Range("D1").Value = "1/19"
